# Updated cryptos list on Tue Jun  6 15:44:47 UTC 2023 with GitHub Actions
#
# Rows 2-37 keep the same coin/link but get refreshed Price (D) and
# Volume(1h) (E) figures. Rows 38-51 gain a new "Frax" entry at row 38,
# pushing RenderToken..Cronos down by one row and dropping the former
# last entry (NEARProtocol) off the bottom of the A1:E51 table.
#
# Price values are plain, unformatted numeric-looking text in this sheet
# (t="inlineStr", e.g. "0.8040", "1.000", "0.05910"). A bare .Value=
# assignment lets Excel auto-convert that text to a Double, silently
# dropping the significant trailing zeros, so the Price column is written
# via a helper that forces the text number format first; every other
# column (never numeric-looking) is assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Rows 2-37: refresh Price / Volume(1h) only ---
Set-PriceText 'D2' '25.983.27'
$ws.Range('E2').Value = '  -2.15%  '
Set-PriceText 'D3' '1.839.89'
$ws.Range('E3').Value = '  -1.14%  '
Set-PriceText 'D4' '0.9997'
$ws.Range('E4').Value = '  -0.33%  '
Set-PriceText 'D5' '278.14'
$ws.Range('E5').Value = '  -3.03%  '
Set-PriceText 'D6' '0.9998'
$ws.Range('E6').Value = '  -0.40%  '
Set-PriceText 'D7' '0.5096'
$ws.Range('E7').Value = '  -2.43%  '
Set-PriceText 'D8' '0.3496'
$ws.Range('E8').Value = '  -5.16%  '
Set-PriceText 'D9' '44.81'
$ws.Range('E9').Value = '  -0.30%  '
Set-PriceText 'D10' '0.06821'
$ws.Range('E10').Value = '  -3.70%  '
Set-PriceText 'D11' '19.88'
$ws.Range('E11').Value = '  -5.49%  '
Set-PriceText 'D12' '0.8040'
$ws.Range('E12').Value = '  -8.16%  '
Set-PriceText 'D13' '0.07773'
$ws.Range('E13').Value = '  -3.24%  '
Set-PriceText 'D14' '1.835.23'
$ws.Range('E14').Value = '  -1.52%  '
Set-PriceText 'D15' '5.074'
$ws.Range('E15').Value = '  -2.91%  '
Set-PriceText 'D16' '88.14'
$ws.Range('E16').Value = '  -2.77%  '
Set-PriceText 'D17' '0.9993'
$ws.Range('E17').Value = '  -0.34%  '
Set-PriceText 'D18' '14.14'
$ws.Range('E18').Value = '  -2.64%  '
Set-PriceText 'D19' '0.000008062'
$ws.Range('E19').Value = '  -3.98%  '
Set-PriceText 'D20' '1.000'
$ws.Range('E20').Value = '  -0.34%  '
Set-PriceText 'D21' '26.026.59'
$ws.Range('E21').Value = '  -2.13%  '
Set-PriceText 'D22' '4.774'
$ws.Range('E22').Value = '  -2.83%  '
Set-PriceText 'D23' '10.06'
$ws.Range('E23').Value = '  -4.31%  '
Set-PriceText 'D24' '6.205'
$ws.Range('E24').Value = '  -1.37%  '
Set-PriceText 'D25' '2.373'
$ws.Range('E25').Value = '  +6.38%  '
Set-PriceText 'D26' '143.51'
$ws.Range('E26').Value = '  -0.73%  '
Set-PriceText 'D27' '1.663'
$ws.Range('E27').Value = '  -4.04%  '
Set-PriceText 'D28' '17.17'
$ws.Range('E28').Value = '  -3.45%  '
Set-PriceText 'D29' '109.58'
$ws.Range('E29').Value = '  -2.79%  '
Set-PriceText 'D30' '4.359'
$ws.Range('E30').Value = '  -5.82%  '
Set-PriceText 'D31' '4.281'
$ws.Range('E31').Value = '  -5.78%  '
Set-PriceText 'D32' '0.08805'
$ws.Range('E32').Value = '  -2.20%  '
Set-PriceText 'D33' '0.04855'
$ws.Range('E33').Value = '  -1.44%  '
Set-PriceText 'D34' '1.161'
$ws.Range('E34').Value = '  +0.92%  '
Set-PriceText 'D35' '0.7264'
$ws.Range('E35').Value = '  -7.20%  '
Set-PriceText 'D36' '2.863'
$ws.Range('E36').Value = '  -1.59%  '
Set-PriceText 'D37' '3.202'
$ws.Range('E37').Value = '  +1.47%  '

# --- Rows 38-51: Frax inserted, remaining coins shift down one row ---
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-PriceText 'D38' '0.9989'
$ws.Range('E38').Value = '  -0.84%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-PriceText 'D39' '2.351'
$ws.Range('E39').Value = '  -9.20%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-PriceText 'D40' '0.01848'
$ws.Range('E40').Value = '  -3.61%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-PriceText 'D41' '0.5144'
$ws.Range('E41').Value = '  -12.22%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-PriceText 'D42' '0.9467'
$ws.Range('E42').Value = '  -8.66%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-PriceText 'D43' '117.05'
$ws.Range('E43').Value = '  +2.45%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-PriceText 'D44' '6.248'
$ws.Range('E44').Value = '  -2.52%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-PriceText 'D45' '8.011'
$ws.Range('E45').Value = '  -6.44%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-PriceText 'D46' '0.9993'
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-PriceText 'D47' '0.4511'
$ws.Range('E47').Value = '  -11.74%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-PriceText 'D48' '0.1360'
$ws.Range('E48').Value = '  -7.43%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-PriceText 'D49' '9.269'
$ws.Range('E49').Value = '  -5.79%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-PriceText 'D50' '36.13'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-PriceText 'D51' '0.05910'
$ws.Range('E51').Value = '  -1.66%  '
